$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on columns D and E so numeric-looking strings (prices, percents)
# are preserved exactly as literal text instead of being parsed into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '27.384.40'
$ws.Range('D3').Value = '1.827.33'
$ws.Range('E3').Value = '  +1.10%  '
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '313.21'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = '0.4621'
$ws.Range('E7').Value = '  +3.88%  '
$ws.Range('D8').Value = '0.3779'
$ws.Range('E8').Value = '  +2.82%  '
$ws.Range('D9').Value = '0.07415'
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('D10').Value = '0.8777'
$ws.Range('E10').Value = '  +2.52%  '
$ws.Range('E11').Value = '  +0.67%  '
$ws.Range('D12').Value = '1.828.49'
$ws.Range('E12').Value = '  +0.74%  '
$ws.Range('D13').Value = '6.722'
$ws.Range('D14').Value = '5.444'
$ws.Range('E14').Value = '  +2.57%  '
$ws.Range('D15').Value = '93.14'
$ws.Range('E15').Value = '  +0.68%  '
$ws.Range('D16').Value = '0.07078'
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').Value = '0.000008804'
$ws.Range('E18').Value = '  +0.82%  '
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('D20').Value = '15.07'
$ws.Range('E20').Value = '  +1.43%  '
$ws.Range('D21').Value = '27.382.06'
$ws.Range('E21').Value = '  +1.73%  '
$ws.Range('D22').Value = '5.336'
$ws.Range('E22').Value = '  +3.57%  '
$ws.Range('D23').Value = '10.96'
$ws.Range('E23').Value = '  +1.29%  '
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D24').Value = '2.053.39'
$ws.Range('E24').Value = '  -2.45%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '1.953'
$ws.Range('E25').Value = '  -2.11%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '151.18'
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '2.266'
$ws.Range('E27').Value = '  +3.99%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '18.59'
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '5.349'
$ws.Range('E29').Value = '  +2.66%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '117.12'
$ws.Range('E30').Value = '  +0.50%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.08949'
$ws.Range('E31').Value = '  +1.38%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '0.8010'
$ws.Range('E32').Value = '  +6.86%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').Value = '1.196'
$ws.Range('E33').Value = '  +1.89%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '4.552'
$ws.Range('E34').Value = '  +2.09%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '2.934'
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').Value = '0.9996'
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '1.099'
$ws.Range('E37').Value = '  +1.07%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.01979'
$ws.Range('E38').Value = '  +0.65%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.05267'
$ws.Range('E39').Value = '  +1.51%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '7.366'
$ws.Range('E40').Value = '  +5.05%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '0.5343'
$ws.Range('E41').Value = '  +0.44%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '2.380'
$ws.Range('E42').Value = '  +19.66%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '2.892'
$ws.Range('E43').Value = '  +0.89%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').Value = '0.1707'
$ws.Range('E44').Value = '  +0.95%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').Value = '8.673'
$ws.Range('E45').Value = '  +3.04%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '0.5112'
$ws.Range('E46').Value = '  -0.90%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '10.58'
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '105.51'
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '1.685'
$ws.Range('E49').Value = '  +1.07%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').Value = '0.9993'
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.06379'
$ws.Range('E51').Value = '  +3.88%  '
